$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells that look numeric stay literal text (matches source formatting,
# e.g. "27.717.04" / "319.70" / "0.3746" are text, not numbers).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.717.04"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.848.63"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  -2.71%  "
$ws.Range("D5").Value = "319.70"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("E6").Value = "  -2.43%  "
$ws.Range("E7").Value = "  -2.45%  "
$ws.Range("D8").Value = "0.3746"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("D10").Value = "0.8803"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").Value = "21.59"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "1.851.76"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "6.732"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "5.454"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "0.07137"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").Value = "87.94"
$ws.Range("E16").Value = "  +4.93%  "
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("D18").Value = "0.000008991"
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("D20").Value = "15.49"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "27.736.48"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").Value = "11.16"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").Value = "2.076.54"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").Value = "2.015"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").Value = "155.77"
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("D27").Value = "18.57"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("D28").Value = "2.134"
$ws.Range("E28").Value = "  +7.58%  "
$ws.Range("D29").Value = "5.391"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").Value = "120.42"
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").Value = "0.7792"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "4.562"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "2.917"
$ws.Range("E35").Value = "  -6.09%  "
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("D37").Value = "1.140"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").Value = "0.05334"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").Value = "0.01975"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Value = "7.281"
$ws.Range("E40").Value = "  +5.62%  "
$ws.Range("D41").Value = "2.875"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "0.5164"
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").Value = "8.905"
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("D47").Value = "0.4732"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "0.06510"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").Value = "1.699"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").Value = "1.011"
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("D51").Value = "1.881"
$ws.Range("E51").Value = "  -2.04%  "

# Rows 45/46 swap identity: Quant moves up to row 45, EnergySwap moves to row 46,
# each carrying updated price/volume figures.
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "109.38"
$ws.Range("E45").Value = "  -0.58%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "10.65"
$ws.Range("E46").Value = "  -0.61%  "
